$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Link the keynote speaker (Julia Lane) to her new DIFA project speaker page
# instead of her NYU Wagner faculty page.
$ws.Range("F4").Value = "[Julia Lane](https://dataifa.github.io/difa-project/julia_lane.html)"

# Update the active selection to match the saved view.
$ws.Range("F14").Select()
